$d = $word.ActiveDocument

# 1. Fix "quis" -> "quiz" in the dashboard/quiz sentence
$d.Content.Find.Execute(
    "através de acertos em um quis, além de possibilitar",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "através de acertos em um quiz, além de possibilitar", 2)

# 2. Bold the GitHub repository hyperlink text
$r = $d.Content
$r.Find.Execute("https://github.com/bruno-yuji/Projeto_Individual", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Bold = 1
